$d = $word.ActiveDocument
$sel = $word.Selection

# Replace the text of paragraphs 1-7 (1-based) in place, preserving the
# paragraph mark, using Selection.SetRange + TypeText so that straight
# quotes/apostrophes are not auto-corrected to curly quotes and so that
# xml:space="preserve" is (re)computed from the final text.

$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$sel.SetRange($r1.Start, $r1.End - 1) | Out-Null
$sel.TypeText('⚡️🚀המאמר היומי של מייק 14.08.24: ⚡️🚀')

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$sel.SetRange($r2.Start, $r2.End - 1) | Out-Null
$sel.TypeText('Jumping Ahead: Improving Reconstruction Fidelity with JumpReLU Sparse Autoencoders')

$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$sel.SetRange($r3.Start, $r3.End - 1) | Out-Null
$sel.TypeText('אתמול סקרנו מאמר שהשתמש בגישת SAE או Sparse AutoEncoders כדי לחדור ל״מחשובותיו״ של מודל שפה גדול דרך האקטיבציות של הנוירונים שלהם. הנחת היסוד במאמר היתה כי נוירונים ״מגיבים״ לכמה קונספטים שונים וניתן לאמן SAE רדוד מאוד (שכבה אחת בדקודר ושכבה אחת באנדוקר) כדי להגיע לוקטור דליל המקודד (נדלק) קונספט אחד בלבד כלומר disentanglement של הפיצ''רים לנוירונים ייעודיים.')

$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$sel.SetRange($r4.Start, $r4.End - 1) | Out-Null
$sel.TypeText('כמאמר יש באנקודר של SAE שכבה לינארית אחת עם פונקציית אקטיבציה הנקראת JumpReLU שראיתי אותה בפעם הראשונה במאמר הזה. פונקציה הזו היא בעצם הזזה של ReLU בציר X ובציר y בפרמטר t נלמד (במאמר זה נקרא טטה). הטענה במאמר שזה מאפשר ללמוד את הייצוג הדליל של דאטה על ידי האנקודר יותר טוב של פונקציית ReLU בגלל שהוא מאפשר לאפס את הקטיבציות בצורה ״נלמדת יותר מ-ReLU".')

$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$sel.SetRange($r5.Start, $r5.End - 1) | Out-Null
$sel.TypeText('עכשיו נשאלת השאלה איך אנחנו אוכפים דלילות על ייצוג הדאטה (אחרי האנקודר). בעבודות קודמות השתמשו ב-L1 בשביל כך אך כאן המחברים משתמשים באותה JumpReLU כדי להפוך את איפוס האיברים בייצוג יותר נלמד. ושימו לב ש- JumpReLU בא עם פרמטר נלמד הזה לזה של האנקודר עצמו שזה עוזר לאכוף דלילות על הייצוג.')

$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$sel.SetRange($r6.Start, $r6.End - 1) | Out-Null
$sel.TypeText('יש עוד טריק אחד קטן ולא מאוד מהותי במאמר הנקרא Kernel density estimation או KDE. אם אתם זוכרים KDE עוזר לנו לשערך(כלומר לקרב) פונקצית צפיפות בהינתם דאטהסט של נקודות באמצעות פונקציית קרנל. פונקציית קרנל יכולה להיות גאוסית למשל ומטרתה לשערך את פונקציית הצפיפות לנקודות לא ידועות על ידי קירובה בין הנקודות בדאטהסט (בדומה לספליין). אז המחברים משתמשים בטריק הזה כדי לשערך את JumpReLU בנקודה t שבה היא לא גזירה. ')

$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
$sel.SetRange($r7.Start, $r7.End - 1) | Out-Null
$sel.TypeText('מאמר נחמד בנושא די חשוב שאמשיך לסקור כנראה גם בעתיד…')

# Append a brand new 8th paragraph holding the new link, after the
# (7th, now-last) paragraph.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newLastPara = $d.Paragraphs.Last
$sel.SetRange($newLastPara.Range.Start, $newLastPara.Range.End - 1) | Out-Null
$sel.TypeText('https://arxiv.org/pdf/2407.14435')

Write-Output "Edit complete"
